$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the per-observation data (Id, Antal, Alder-Stadium,
# Ost, Nord, Publik kommentar) that is currently on row 2 with the data
# currently on row 4, while rows 1 and 3 stay untouched.
#
# We perform the swap with Copy (instead of direct .Value2 assignment)
# so that cell data types (e.g. text-formatted numbers like "25"/"500")
# and blank cells are preserved faithfully, using a scratch cell far
# outside the used range as temporary storage for each 3-way swap.

$columns = @("A", "I", "K", "Q", "R", "AC")
$tempCell = "ZZ100"

foreach ($col in $columns) {
    $src = $ws.Range($col + "2")
    $dst = $ws.Range($col + "4")
    $tmp = $ws.Range($tempCell)

    # Remember whether the source cell (row 2) was actually blank, since
    # Copy()-ing a blank cell onto a non-blank one is a no-op in this
    # engine (it will not overwrite/clear the destination's content).
    $srcWasBlank = ($src.Value2 -eq $null) -or ($src.Value2 -eq "")

    $src.Copy($tmp)
    $dst.Copy($src)
    $tmp.Copy($dst)
    $tmp.ClearContents()

    if ($srcWasBlank) {
        $dst.ClearContents()
    }
}
